$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (match data) with the new fixture/odds values
$ws.Range("A2").Value = "IBJAWlHc"
$ws.Range("B2").Value = "25/11/2024"
$ws.Range("C2").Value = "11:00"
$ws.Range("D2").Value = "INDIA - ISL"
$ws.Range("E2").Value = "Hyderabad"
$ws.Range("F2").Value = "Odisha FC"
$ws.Range("G2").Value = 2.75
$ws.Range("H2").Value = 3.75
$ws.Range("I2").Value = 2.05
$ws.Range("J2").Value = 3.6
$ws.Range("K2").Value = 2.25
$ws.Range("L2").Value = 2.75
$ws.Range("M2").Value = 1.03
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.08
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.67
$ws.Range("V2").Value = 2.1
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 15
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 7.5
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 11
$ws.Range("AJ2").Value = 9
$ws.Range("AK2").Value = 21
$ws.Range("AL2").Value = 17
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 5
$ws.Range("AO2").Value = 17
$ws.Range("AP2").Value = 23
$ws.Range("AQ2").Value = 51
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 3.25
$ws.Range("AU2").Value = 7.5
$ws.Range("AV2").Value = 51
$ws.Range("AW2").Value = 4.33
$ws.Range("AX2").Value = 12
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 41
$ws.Range("BA2").Value = 51
$ws.Range("BB2").Value = 126

# Remove the trailing "Odd_CS_3-3_HT" / "Odd_CS_4-4_HT" columns (BC:BD)
# entirely -- this clears their header + data cells and shifts the
# sheet dimension back down to column BB.
$ws.Columns("BC:BD").Delete()
